$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "Country Auth."
$ws.Range("B14").Value = "Country Authority"
$ws.Range("C14").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(B14)," ","_"),"-","_")'
$ws.Range("A29").Value = "Edit Status"
$ws.Range("B29").Value = "Edit Status"
$ws.Range("C29").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(B29)," ","_"),"-","_")'
$ws.Range("C29").Select() | Out-Null
